# Update countries & provincias Spain
# Refresh the COVID dashboard numbers for several countries and re-sort the
# "Israel / Barein / Rumania / Nigeria / Armenia" block (Armenia's updated
# total now outranks Rumania/Nigeria, so it moves up the table). Also swap
# the Dominica/Fiyi labels (tied case counts, order flipped) and bump the
# "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Estados Unidos (row 4): totals refreshed ---------------------------
$ws.Cells.Item(4, 2).Value = 2935982
$ws.Cells.Item(4, 3).Value = 212
$ws.Cells.Item(4, 4).Value = 1260472
$ws.Cells.Item(4, 5).Value = 1543192

# --- India (row 7): totals refreshed -------------------------------------
$ws.Cells.Item(7, 5).Value = 245940
$ws.Cells.Item(7, 7).Value = 10
$ws.Cells.Item(7, 8).Value = 19289

# --- Singapur (row 40): totals refreshed ---------------------------------
$ws.Cells.Item(40, 2).Value = 44800
$ws.Cells.Item(40, 3).Value = 136
$ws.Cells.Item(40, 5).Value = 4657

# --- Israel (row 49): totals refreshed -----------------------------------
$ws.Cells.Item(49, 2).Value = 29366
$ws.Cells.Item(49, 3).Value = 196
$ws.Cells.Item(49, 4).Value = 17847
$ws.Cells.Item(49, 5).Value = 11189

# --- Rows 51-53: Armenia's refreshed total (28606) now beats Rumania
# (28582) and Nigeria (28167), so it moves to row 51; Rumania and Nigeria
# each shift down one row, keeping their own (unchanged) figures. --------
$rumaniaRow = @($ws.Cells.Item(51, 2).Value(), $ws.Cells.Item(51, 3).Value(), $ws.Cells.Item(51, 4).Value(), $ws.Cells.Item(51, 5).Value(), $ws.Cells.Item(51, 6).Value(), $ws.Cells.Item(51, 7).Value(), $ws.Cells.Item(51, 8).Value())
$nigeriaRow = @($ws.Cells.Item(52, 2).Value(), $ws.Cells.Item(52, 3).Value(), $ws.Cells.Item(52, 4).Value(), $ws.Cells.Item(52, 5).Value(), $ws.Cells.Item(52, 6).Value(), $ws.Cells.Item(52, 7).Value(), $ws.Cells.Item(52, 8).Value())

$ws.Cells.Item(51, 1).Value = "Armenia"
$ws.Cells.Item(51, 2).Value = 28606
$ws.Cells.Item(51, 3).Value = 706
$ws.Cells.Item(51, 4).Value = 16140
$ws.Cells.Item(51, 5).Value = 11982
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 7
$ws.Cells.Item(51, 8).Value = 484

$ws.Cells.Item(52, 1).Value = "Rumania"
$ws.Cells.Item(52, 2).Value = $rumaniaRow[0]
$ws.Cells.Item(52, 3).Value = $rumaniaRow[1]
$ws.Cells.Item(52, 4).Value = $rumaniaRow[2]
$ws.Cells.Item(52, 5).Value = $rumaniaRow[3]
$ws.Cells.Item(52, 6).Value = $rumaniaRow[4]
$ws.Cells.Item(52, 7).Value = $rumaniaRow[5]
$ws.Cells.Item(52, 8).Value = $rumaniaRow[6]

$ws.Cells.Item(53, 1).Value = "Nigeria"
$ws.Cells.Item(53, 2).Value = $nigeriaRow[0]
$ws.Cells.Item(53, 3).Value = $nigeriaRow[1]
$ws.Cells.Item(53, 4).Value = $nigeriaRow[2]
$ws.Cells.Item(53, 5).Value = $nigeriaRow[3]
$ws.Cells.Item(53, 6).Value = $nigeriaRow[4]
$ws.Cells.Item(53, 7).Value = $nigeriaRow[5]
$ws.Cells.Item(53, 8).Value = $nigeriaRow[6]

# --- Rows 205-206: Dominica/Fiyi are tied (18 cases each) and swap order -
$ws.Cells.Item(205, 1).Value = "Fiyi"
$ws.Cells.Item(206, 1).Value = "Dominica"

# --- Footer timestamp (row 1) --------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 5 de Julio de 2020 a las 10:19"
